$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated coin data (price + 1h volume change) scraped on
# Sun Mar 19 03:27:04 UTC 2023.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.216.75"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.54%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.784.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.70%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.19%  "

$ws.Range("E6").Value = "  -0.03%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3778"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.21%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3434"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.16%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.30"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.96%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.196"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.86%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07472"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.36%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.14%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.69"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.99%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.459"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.39%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.789.26"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.25%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.085"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.71%  "

$ws.Range("E17").Value = "  -3.61%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06664"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.75%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "83.89"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.98%  "

$ws.Range("E20").Value = "  -0.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.619"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.22%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "17.28"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.38%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "27.241.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.44%  "

$ws.Range("E24").Value = "  -5.91%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.421"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.06%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.508"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.86%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.539"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -6.97%  "

$ws.Range("E28").Value = "  -3.85%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "153.62"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.31%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.989.60"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.35%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "133.98"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.97%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.010"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.24%  "

$ws.Range("E33").Value = "  -6.99%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08679"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.14%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "13.22"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.39%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.658"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.86%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6940"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.40%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.456"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.46%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06328"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -4.11%  "

$ws.Range("E40").Value = "  -3.58%  "

$ws.Range("B41").Value = "VeChain"
$ws.Range("C41").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.02339"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.10%  "

$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.778"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.49%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.243"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.36%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.45"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.21%  "

$ws.Range("E45").Value = "  -1.83%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.000"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.14%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.853"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.38%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.144"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.63%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "129.06"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.54%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07125"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.06%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "79.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.46%  "
